# Update genotype labels in B49:B59 for readability (final sheet for SFN tweet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 49-53 were "E1 HET;E2cKO" -> now "E1 HET; E2 cKO"
foreach ($r in 49..53) {
    $ws.Cells.Item($r, 2).Value = "E1 HET; E2 cKO"
}

# Rows 54-59 were "E1ko; E2 HET" -> now "E1 ko; E2 HET"
foreach ($r in 54..59) {
    $ws.Cells.Item($r, 2).Value = "E1 ko; E2 HET"
}

# Update the active selection to B62
$ws.Range("B62").Select()
